# Update the marksheet "Corr/total marks" figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking row: correct-answer count goes from 3 to 5
$ws.Range("B11").Value = 5

# Total row: total marks go from 57 to 95
$ws.Range("B12").Value = 95

# Total row: corr/total text goes from "53/84" to "95/140"
$ws.Range("E12").Value = "95/140"
